# Weekly Fruta/Hortaliza update: insert two new daily price rows for
# "Arándano (blue)" at Mercado Mayorista Lo Valledor de Santiago
# (week of 2022-12-23), pushing the existing rows 443:490 down to 445:492.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 443 so the existing data
# (previously rows 443:490) shifts down to rows 445:492.
$ws.Rows("443:444").Insert()

# --- New row 443 ---
$ws.Range("A443").Value = 6
$ws.Range("B443").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C443").Value = "Metropolitana"
$ws.Range("D443").Value = 44918
$ws.Range("E443").Value = 13
$ws.Range("F443").Value = "Fruta"
$ws.Range("G443").Value = 100101
$ws.Range("H443").Value = "Berries"
$ws.Range("I443").Value = 100101001
$ws.Range("J443").Value = "Arándano (blue)"
$ws.Range("K443").Value = "Sin especificar"
$ws.Range("L443").Value = "Especial"
$ws.Range("M443").Value = 500
$ws.Range("N443").Value = 3000
$ws.Range("O443").Value = 3000
$ws.Range("P443").Value = 3000
$ws.Range("Q443").Value = "`$/bandeja 2 kilos"
$ws.Range("R443").Value = "Provincia de Curicó"
$ws.Range("S443").Value = 1500
$ws.Range("T443").Value = 2

# --- New row 444 ---
$ws.Range("A444").Value = 6
$ws.Range("B444").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C444").Value = "Metropolitana"
$ws.Range("D444").Value = 44918
$ws.Range("E444").Value = 13
$ws.Range("F444").Value = "Fruta"
$ws.Range("G444").Value = 100101
$ws.Range("H444").Value = "Berries"
$ws.Range("I444").Value = 100101001
$ws.Range("J444").Value = "Arándano (blue)"
$ws.Range("K444").Value = "Sin especificar"
$ws.Range("L444").Value = "Especial"
$ws.Range("M444").Value = 750
$ws.Range("N444").Value = 3000
$ws.Range("O444").Value = 3000
$ws.Range("P444").Value = 3000
$ws.Range("Q444").Value = "`$/bandeja 2 kilos"
$ws.Range("R444").Value = "Región del Maule"
$ws.Range("S444").Value = 1500
$ws.Range("T444").Value = 2
